$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 54, pushing the existing rows 54..78 down to 55..79.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly entry.
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 44777
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = 100114001
$ws.Range("G54").Value = "Papa"
$ws.Range("H54").Value = "Asterix"
$ws.Range("I54").Value = "1a (guarda)"
$ws.Range("J54").Value = 1000
$ws.Range("K54").Value = 11000
$ws.Range("L54").Value = 12000
$ws.Range("M54").Value = 11500
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Región de Los Lagos"
$ws.Range("P54").Value = 460
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
